$d = $word.ActiveDocument

# Locate the paragraph that currently holds the {{CONTENT}} field group
# immediately followed by a {{TITLE}} field group inside the very same
# paragraph (no paragraph break between them). That paragraph is the one
# that needs a TITLE-only paragraph and a fresh CONTENT field group
# inserted in front of it.
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $codes = @()
    for ($j = 1; $j -le $d.Fields.Count; $j++) {
        $fld = $d.Fields.Item($j)
        if ($fld.Code.Start -ge $para.Range.Start -and $fld.Code.End -le $para.Range.End) {
            $codes += $fld.Code.Text.Trim()
        }
    }
    if ($codes.Count -ge 2 -and $codes[0] -eq "{{CONTENT}}" -and $codes[1] -eq "{{TITLE}}") {
        $targetParaIndex = $i
        break
    }
}

if ($targetParaIndex -eq -1) {
    # Fallback: the paragraph that contains more than one field is the one
    # that needs splitting.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $count = 0
        for ($j = 1; $j -le $d.Fields.Count; $j++) {
            $fld = $d.Fields.Item($j)
            if ($fld.Code.Start -ge $para.Range.Start -and $fld.Code.End -le $para.Range.End) {
                $count += 1
            }
        }
        if ($count -ge 2) {
            $targetParaIndex = $i
            break
        }
    }
}

$targetPara = $d.Paragraphs.Item($targetParaIndex)
$insertPos = $targetPara.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

# Insert a standalone "${TITLE}" field-code paragraph, immediately followed
# by a standalone "${CONTENT}" field-code paragraph whose closing mark is
# left open so it merges with (and now leads) the paragraph's pre-existing
# {{CONTENT}} / {{TITLE}} field groups.
$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:fldChar w:fldCharType="begin"/>
            </w:r>
            <w:r>
              <w:instrText xml:space="preserve"> ${TITLE} </w:instrText>
            </w:r>
            <w:r>
              <w:fldChar w:fldCharType="separate"/>
            </w:r>
            <w:r>
              <w:fldChar w:fldCharType="end"/>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:fldChar w:fldCharType="begin"/>
            </w:r>
            <w:r>
              <w:instrText xml:space="preserve"> ${CONTENT} </w:instrText>
            </w:r>
            <w:r>
              <w:fldChar w:fldCharType="separate"/>
            </w:r>
            <w:r>
              <w:fldChar w:fldCharType="end"/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xmlFragment)
